$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 (spreadsheet row 15): Глава 14, "3114-245", Котел,
# "Мегу, Сира, Эван, Барри, Люк, Зак, Джун, Широ",
# "Экзамен в захватчики, бой с Эваном"
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "3114-245"
$ws.Range("C15").Value = "Котел"
$ws.Range("D15").Value = "Мегу, Сира, Эван, Барри, Люк, Зак, Джун, Широ"
$ws.Range("E15").Value = "Экзамен в захватчики, бой с Эваном"

# Match formatting used by the row above it (row 14): vertically centered,
# wrapped text cells, row height 30.
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("C15").VerticalAlignment = -4108
$ws.Range("C15").WrapText = $true
$ws.Range("D15").VerticalAlignment = -4108
$ws.Range("D15").WrapText = $true
$ws.Range("E15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 30

# Update selection to reflect where the user ended up after entering the row
$null = $ws.Range("E16").Select()

Write-Output "done"
